# start of season update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manager_ids")

# Fill in ids for managers newly added in the "Season 24-25" block (rows 40-46)
$ws.Range("A40").Value = 1114016
$ws.Range("A41").Value = 109884
$ws.Range("A42").Value = 6188948
$ws.Range("A43").Value = 1123161
$ws.Range("A44").Value = 10946
$ws.Range("A45").Value = 6186573
$ws.Range("A46").Value = 1077612

# Replace the placeholder "-" league_id values with the real league id for the new season
$ws.Range("D38:D46").Value = 20232

# Match the style of column D to the rest of the table (same style as column A-C)
$ws.Range("D1:D37").HorizontalAlignment = -4131
